$wb = $excel.ActiveWorkbook

# --- FT sheet: bump the Final Time year from 2030 to 2050 ---
$wsFT = $wb.Worksheets.Item("FT")
$wsFT.Range("B2").Value = 2050

# The FT sheet is no longer the tab the file opens on -- reset its lingering
# selection away from B3 before handing the active tab back to "About".
$wsFT.Range("A1").Select() | Out-Null

# --- "About" becomes the active/selected sheet again (was "FT") ---
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Activate() | Out-Null

# --- Workbook recalculation settings: switch to manual with iterative calc ---
$excel.Calculation = -4135          # xlCalculationManual
$excel.Iteration = $true
$excel.MaxIterations = 100
$excel.MaxChange = 0.00001
